$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Query 2 (table row 3): mark the "faturacaoGlobal" column (col 4) with an X
$cellFaturacao = $t.Cell(3, 4)
$cellFaturacao.Range.InsertBefore("X")

# Query 4 (table row 5): mark the "Vendas/Filial" column (col 5) with an X,
# and move the "_GoBack" bookmark onto this cell's paragraph (it previously
# sat in the trailing empty paragraph after the table).
$cellVendas = $t.Cell(5, 5)
$rngVendas = $cellVendas.Range
$d.Bookmarks.Add("_GoBack", $rngVendas)
$t.Cell(5, 5).Range.InsertBefore("X")
